$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("D16").Style = "Normal"

$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "642,530,686,576"
$ws.Range("I16").Style = "Normal"

$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.75"
$ws.Range("J16").Style = "Normal"

# Row 17 updates
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("D17").Style = "Normal"

$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "794,481,831,526"
$ws.Range("I17").Style = "Normal"
